$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------
# Grab the formatting needed for the new row 52 first, while the
# "applyFill" style variants (13/14) used by it still exist on
# rows 45-50 in their original, unedited form.
# ---------------------------------------------------------------
CopyFormat "C2"  "C52"
CopyFormat "F45" "D52"
CopyFormat "C45" "E52"
CopyFormat "F45" "F52"
CopyFormat "C45" "G52"
CopyFormat "F45" "H52"
CopyFormat "C45" "I52"

# ---------------------------------------------------------------
# Re-normalize the formatting of the Chapter 9 block (rows 44-51)
# so that it reuses the same style entries as the other chapter
# blocks (e.g. rows 15-21), instead of the chapter-9-only style
# duplicates. This removes the need for the extra cellXfs entries.
# ---------------------------------------------------------------

# Row 44 - first row of the chapter block
CopyFormat "B15" "B44"
CopyFormat "C15" "C44"
CopyFormat "D15" "D44"
CopyFormat "E15" "E44"
CopyFormat "F15" "F44"
CopyFormat "G15" "G44"

# Rows 45-48 - plain interior rows
foreach ($r in 45..48) {
    CopyFormat "C16" "C$r"
    CopyFormat "D1"  "D$r"
    CopyFormat "E16" "E$r"
    CopyFormat "F1"  "F$r"
    CopyFormat "G16" "G$r"
    CopyFormat "H1"  "H$r"
    $ws.Range("I$r`:J$r").Clear()
}

# Row 49 - interior row that keeps a Function value in column I
CopyFormat "C16" "C49"
CopyFormat "D1"  "D49"
CopyFormat "E16" "E49"
CopyFormat "F1"  "F49"
CopyFormat "G16" "G49"
CopyFormat "H1"  "H49"
CopyFormat "C16" "I49"
$ws.Range("J49").Clear()

# Row 50 - interior row, no Function/Description values
CopyFormat "C16" "C50"
CopyFormat "D1"  "D50"
CopyFormat "E16" "E50"
CopyFormat "F1"  "F50"
CopyFormat "G16" "G50"
CopyFormat "H1"  "H50"
$ws.Range("I50:J50").Clear()

# Row 51 - last row of the chapter block
CopyFormat "C21" "C51"
CopyFormat "E21" "E51"
CopyFormat "F21" "F51"
CopyFormat "G21" "G51"
CopyFormat "H21" "H51"

# ---------------------------------------------------------------
# Add new row 52: Chapter 10, Section 1, Subsection 1
# ---------------------------------------------------------------
$ws.Range("B52").Formula = '=_xlfn.CONCAT(TEXT(C52,"00"),TEXT(E52,"00"),TEXT(G52,"00"))'
$ws.Range("C52").Value = 10
$ws.Range("D52").Value = "Working with the filter context"
$ws.Range("E52").Value = 1
$ws.Range("F52").Value = "Using HASONEVALUE and SELECTEDVALUE"
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = "Using HASONEVALUE and SELECTEDVALUE"
$ws.Range("I52").Value = "HASONEVALUE, SELECTEDVALUE"

# Update the selected cell to match the new selection in the sheet view
$ws.Range("H58").Select()
